$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-02 Monday" "2024-12-03 Tuesday"

Replace-Text "784×7=" "200×7="
Replace-Text "892×5=" "515×5="
Replace-Text "752×8=" "487×9="
Replace-Text "431×9=" "879×4="
Replace-Text "489×2=" "276×2="

Replace-Text "610×6=" "947×8="
Replace-Text "257×4=" "624×2="
Replace-Text "445×6=" "378×7="
Replace-Text "427×2=" "811×9="
Replace-Text "131×9=" "795×3="

Replace-Text "984×4=" "708×4="
Replace-Text "209×2=" "249×4="
Replace-Text "714×2=" "210×5="
Replace-Text "665×6=" "177×7="
Replace-Text "700×2=" "897×9="

Replace-Text "175×2=" "135×2="
Replace-Text "911×3=" "887×6="
Replace-Text "650×5=" "369×5="
Replace-Text "359×8=" "192×7="
Replace-Text "805×8=" "186×8="

Replace-Text "923×6=" "635×8="
Replace-Text "796×7=" "484×8="
Replace-Text "722×9=" "665×5="
Replace-Text "513×9=" "367×5="
Replace-Text "902×5=" "372×7="
